$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.447.81"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").Value = "3.388.62"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'581.29"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").Value = "'179.06"
$ws.Range("E6").Value = "  +0.89%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.594"
$ws.Range("E8").Value = "  +0.61%  "

$ws.Range("D9").Value = "'0.198"
$ws.Range("E9").Value = "  +7.94%  "

$ws.Range("D10").Value = "'0.587"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("D11").Value = "'48.41"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "'0.0000283"
$ws.Range("E12").Value = "  +3.65%  "

$ws.Range("D13").Value = "'686.51"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "'8.59"
$ws.Range("E14").Value = "  +2.02%  "

$ws.Range("D15").Value = "3.921.89"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").Value = "69.526.55"
$ws.Range("E16").Value = "  +1.77%  "

$ws.Range("D17").Value = "'0.121"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "3.379.13"
$ws.Range("E18").Value = "  +1.09%  "

$ws.Range("D19").Value = "'17.71"
$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("D20").Value = "'11.28"
$ws.Range("E20").Value = "  +0.73%  "

$ws.Range("D21").Value = "'0.909"
$ws.Range("E21").Value = "  +1.61%  "

$ws.Range("D22").Value = "'17.18"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").Value = "'5.35"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").Value = "'101.17"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").Value = "'3.88"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").Value = "'2.70"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'9.73"
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("D28").Value = "'33.51"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").Value = "'8.73"
$ws.Range("E29").Value = "  +2.74%  "

$ws.Range("D30").Value = "'6.95"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").Value = "'3.87"
$ws.Range("E31").Value = "  +17.75%  "

$ws.Range("D32").Value = "'11.03"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("D33").Value = "'549.47"
$ws.Range("E33").Value = "  -2.13%  "

$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").Value = "'57.88"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Value = "3.605.32"
$ws.Range("E37").Value = "  -2.52%  "

$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  +3.28%  "

$ws.Range("D39").Value = "'35.46"
$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("D40").Value = "0.0₃0744"
$ws.Range("E40").Value = "  +10.60%  "

$ws.Range("D41").Value = "'3.33"
$ws.Range("E41").Value = "  +5.18%  "

$ws.Range("D42").Value = "'2.72"
$ws.Range("E42").Value = "  +4.25%  "

$ws.Range("E43").Value = "  +3.58%  "

$ws.Range("D44").Value = "'0.0425"
$ws.Range("E44").Value = "  +3.41%  "

$ws.Range("D45").Value = "'0.336"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "'2.67"
$ws.Range("E46").Value = "  +0.84%  "

$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").Value = "'1.39"
$ws.Range("E48").Value = "  +3.75%  "

$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").Value = "'129.81"
$ws.Range("E50").Value = "  -0.78%  "

$ws.Range("D51").Value = "'2.59"
$ws.Range("E51").Value = "  +1.17%  "

Write-Host "Updated cryptos list"